$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 / Row 3 swap of A/Q/R, plus B column update on rows 2-4 ---
$ws.Range("A2").Value = 112095428
$ws.Range("B2").Value = 78699
$ws.Range("Q2").Value = 491096
$ws.Range("R2").Value = 6954259

$ws.Range("A3").Value = 112095298
$ws.Range("B3").Value = 78699
$ws.Range("Q3").Value = 491104
$ws.Range("R3").Value = 6954282

$ws.Range("B4").Value = 78699

# --- New row 5 ---
$ws.Range("A5").Value = 112426767
$ws.Range("B5").Value = 78699
$ws.Range("C5").Value = "Ovaliderad"
$ws.Range("D5").Value = "NT"
$ws.Range("E5").Value = 6458
$ws.Range("F5").Value = "Lunglav"
$ws.Range("G5").Value = "Lobaria pulmonaria"
$ws.Range("H5").Value = "(L.) Hoffm."
$ws.Range("I5").Style = "Normal"
$ws.Range("K5").Style = "Normal"
$ws.Range("P5").Value = "Lill-Öretjärnen (Lill-Öretjärnen), Jmt"
$ws.Range("Q5").Value = 490949
$ws.Range("R5").Value = 6953753
$ws.Range("S5").Value = 1
$ws.Range("T5").Value = "Jämtland"
$ws.Range("U5").Value = "Berg"
$ws.Range("V5").Value = "Jämtland"
$ws.Range("W5").Value = "Hackås"
$ws.Range("Y5").NumberFormat = "@"
$ws.Range("Y5").Value = "2023-09-30"
$ws.Range("Y5").Style = "Normal"
$ws.Range("Z5").Value = "18:45"
$ws.Range("AA5").NumberFormat = "@"
$ws.Range("AA5").Value = "2023-09-30"
$ws.Range("AA5").Style = "Normal"
$ws.Range("AB5").Value = "18:45"
$ws.Range("AD5").Value = $false
$ws.Range("AE5").Value = $false
$ws.Range("AG5").Value = $false
$ws.Range("AT5").Style = "Normal"
$ws.Range("AW5").Value = "Erik Wilhelmsson"
$ws.Range("AX5").Value = "Erik Wilhelmsson"
$ws.Range("AY5").Style = "Normal"

# --- New row 6 ---
$ws.Range("A6").Value = 112426713
$ws.Range("B6").Value = 78699
$ws.Range("C6").Value = "Ovaliderad"
$ws.Range("D6").Value = "NT"
$ws.Range("E6").Value = 6458
$ws.Range("F6").Value = "Lunglav"
$ws.Range("G6").Value = "Lobaria pulmonaria"
$ws.Range("H6").Value = "(L.) Hoffm."
$ws.Range("I6").Style = "Normal"
$ws.Range("K6").Style = "Normal"
$ws.Range("P6").Value = "Lill-Öretjärnen (Lill-Öretjärnen), Jmt"
$ws.Range("Q6").Value = 490958
$ws.Range("R6").Value = 6953733
$ws.Range("S6").Value = 1
$ws.Range("T6").Value = "Jämtland"
$ws.Range("U6").Value = "Berg"
$ws.Range("V6").Value = "Jämtland"
$ws.Range("W6").Value = "Hackås"
$ws.Range("Y6").NumberFormat = "@"
$ws.Range("Y6").Value = "2023-09-30"
$ws.Range("Y6").Style = "Normal"
$ws.Range("Z6").Value = "18:39"
$ws.Range("AA6").NumberFormat = "@"
$ws.Range("AA6").Value = "2023-09-30"
$ws.Range("AA6").Style = "Normal"
$ws.Range("AB6").Value = "18:39"
$ws.Range("AD6").Value = $false
$ws.Range("AE6").Value = $false
$ws.Range("AG6").Value = $false
$ws.Range("AT6").Style = "Normal"
$ws.Range("AW6").Value = "Erik Wilhelmsson"
$ws.Range("AX6").Value = "Erik Wilhelmsson"
$ws.Range("AY6").Style = "Normal"
